# Update HZNP Quarterly Financials worksheet with the two newest quarters.
# Inserts two new columns (D:E) before the existing data, shifting the prior
# quarters two columns to the right (old D -> F, old E -> G, ... old K -> M),
# then fills in the two new quarters' figures and a handful of revised
# figures for the quarter that is now in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank columns before column D; this shifts D:K -> F:M.
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the shifted former column D) into
# the two newly inserted columns so the new data matches the existing
# date / numeric formatting for every row.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# (coordinate, value) pairs: the two new quarters' data (columns D & E),
# plus a small number of corrected figures in the shifted data (mostly
# column H) that differ from a pure shift of the old values.
$data = @(
    @("D7", 43465),
    @("E7", 43373),
    @("D8", 355500),
    @("E8", 325300),
    @("D9", 107100),
    @("E9", 99000),
    @("D10", 248400),
    @("E10", 226300),
    @("D12", 19700),
    @("E12", 21200),
    @("D13", 0),
    @("E13", 0),
    @("D14", 10800),
    @("E14", 1600),
    @("D15", 0),
    @("E15", 0),
    @("D17", 281900),
    @("E17", 271100),
    @("D18", 73600),
    @("E18", 54200),
    @("D20", -700),
    @("E20", 500),
    @("D21", 141900),
    @("E21", 124000),
    @("D22", 29800),
    @("E22", 30400),
    @("D23", 43100),
    @("E23", 24300),
    @("D24", -9400),
    @("E24", -1700),
    @("D25", 0),
    @("E25", 0),
    @("D26", 52600),
    @("E26", 26000),
    @("D27", 52600),
    @("E27", 26000),
    @("D28", 0),
    @("E28", 0),
    @("D29", 37400),
    @("E29", "NA"),
    @("D30", 0),
    @("E30", 0),
    @("D31", 0),
    @("E31", 0),
    @("D32", 700),
    @("E32", -500),
    @("D33", 89900),
    @("E33", 26000),
    @("D34", 0),
    @("E34", 0),
    @("D35", 89900),
    @("E35", 26000),
    @("D38", 43465),
    @("E38", 43373),
    @("D41", 958700),
    @("E41", 807000),
    @("D42", 0),
    @("E42", 0),
    @("D43", 464700),
    @("E43", 391100),
    @("D44", 55300),
    @("E44", 59100),
    @("D45", 69700),
    @("E45", 81900),
    @("D46", 1548400),
    @("E46", 1339200),
    @("D47", 0),
    @("E47", 0),
    @("D48", 20100),
    @("E48", 16600),
    @("D49", 2551700),
    @("E49", 2635900),
    @("D50", 0),
    @("E50", 0),
    @("D51", 0),
    @("E51", 0),
    @("D52", 26200),
    @("E52", 27700),
    @("D53", 0),
    @("E53", 0),
    @("D54", 4146400),
    @("E54", 4019400),
    @("D57", 30300),
    @("E57", 64800),
    @("D58", 0),
    @("E58", 0),
    @("D59", 731600),
    @("E59", 626800),
    @("D60", 761900),
    @("E60", 691600),
    @("D61", 1896700),
    @("E61", 1890800),
    @("D62", 433600),
    @("E62", 520100),
    @("D63", 0),
    @("E63", 0),
    @("D64", 0),
    @("E64", 0),
    @("D65", 0),
    @("E65", 0),
    @("D66", 3092200),
    @("E66", 3102500),
    @("D68", 0),
    @("E68", 0),
    @("D69", 0),
    @("E69", 0),
    @("D70", 0),
    @("E70", 0),
    @("D71", 0),
    @("E71", 0),
    @("D72", -1314700),
    @("E72", -1414900),
    @("D73", 0),
    @("E73", 0),
    @("D74", 0),
    @("E74", 0),
    @("D75", 0),
    @("E75", 0),
    @("D76", 1054200),
    @("E76", 916900),
    @("D77", 0),
    @("E77", 0),
    @("D80", 43465),
    @("E80", 43373),
    @("D81", 89900),
    @("E81", 26000),
    @("D83", 69000),
    @("E83", 69200),
    @("D84", 0),
    @("E84", 0),
    @("D85", 0),
    @("E85", 0),
    @("D86", 0),
    @("E86", 0),
    @("D87", 0),
    @("E87", 0),
    @("D88", 0),
    @("E88", 0),
    @("D89", 108700),
    @("E89", 84900),
    @("D91", -3900),
    @("E91", -100),
    @("D92", 0),
    @("E92", 0),
    @("D93", 0),
    @("E93", 0),
    @("D94", 31100),
    @("E94", 9300),
    @("D96", 0),
    @("E96", 0),
    @("D97", 0),
    @("E97", 0),
    @("D98", 0),
    @("E98", 0),
    @("D99", 0),
    @("E99", 0),
    @("D100", 9500),
    @("E100", 2400),
    @("D101", -700),
    @("E101", 300),
    @("D102", 148700),
    @("E102", 96800),
    @("H9", 142600),
    @("H10", 131600),
    @("H14", 23700),
    @("H17", 341300),
    @("H18", -67100),
    @("H21", 3200),
    @("H23", -98100),
    @("H26", -112500),
    @("H27", -112500),
    @("H33", -37500),
    @("H35", -37500),
    @("H49", 2874200),
    @("H54", 4202300),
    @("H62", 515000),
    @("H66", 3201000),
    @("H72", -1242100),
    @("H76", 1001300),
    @("H81", -37500),
    @("H89", 143300),
    @("I89", 68300),
    @("I91", -1400),
    @("J91", -1200)
)

foreach ($item in $data) {
    $ws.Range($item[0]).Value = $item[1]
}
